# Update countries & provincias Spain
# - Refresh case totals for several countries (data refresh at 19:05)
# - "Republica del Chad" now outranks "Sierra Leona" et al. (reordered block rows 131-137)
# - "Nepal" now outranks "Liberia" et al. (reordered block rows 142-145)
# - Update the "Datos actualizados..." timestamp banner in A1

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp banner (A1)
$ws.Range("A1").Value = "Datos actualizados a 12 de Mayo de 2020 a las 19:05"

# Row 4: Estados Unidos - updated totals
$ws.Range("B4").Value = 1392976
$ws.Range("C4").Value = 7142
$ws.Range("D4").Value = 263641
$ws.Range("E4").Value = 1047098

# Row 11: Brasil - updated totals
$ws.Range("B11").Value = 172243
$ws.Range("C11").Value = 3100
$ws.Range("E11").Value = 92879
$ws.Range("G11").Value = 355
$ws.Range("H11").Value = 11980

# Row 12: Turquia - updated totals
$ws.Range("B12").Value = 141475
$ws.Range("C12").Value = 1704
$ws.Range("D12").Value = 98889
$ws.Range("E12").Value = 38692
$ws.Range("F12").Value = 1045
$ws.Range("G12").Value = 53
$ws.Range("H12").Value = 3894

# Row 51: Chequia - updated totals
$ws.Range("B51").Value = 8198
$ws.Range("C51").Value = 22
$ws.Range("D51").Value = 4865
$ws.Range("E51").Value = 3050

# Row 55: Marruecos - updated totals
$ws.Range("B55").Value = 6418
$ws.Range("C55").Value = 137
$ws.Range("D55").Value = 2991
$ws.Range("E55").Value = 3239

# Row 104: Sri Lanka - updated totals
$ws.Range("B104").Value = 884
$ws.Range("C104").Value = 21
$ws.Range("E104").Value = 509

# Rows 131-137: "Republica del Chad" jumps ahead of Sierra Leona/Congo/Mauricio/
# Isla de Man/Benin/Montenegro, pushing each of those down one rank.
$ws.Range("A131").Value = "Republica del Chad"
$ws.Range("B131").Value = 357
$ws.Range("C131").Value = 35
$ws.Range("D131").Value = 76
$ws.Range("E131").Value = 241
$ws.Range("F131").Value = 0
$ws.Range("G131").Value = 9
$ws.Range("H131").Value = 40

$ws.Range("A132").Value = "Sierra Leona"
$ws.Range("B132").Value = 338
$ws.Range("C132").Value = 0
$ws.Range("D132").Value = 72
$ws.Range("E132").Value = 247
$ws.Range("F132").Value = 0
$ws.Range("G132").Value = 0
$ws.Range("H132").Value = 19

$ws.Range("A133").Value = "Congo"
$ws.Range("B133").Value = 333
$ws.Range("C133").Value = 0
$ws.Range("D133").Value = 53
$ws.Range("E133").Value = 269
$ws.Range("F133").Value = 0
$ws.Range("G133").Value = 0
$ws.Range("H133").Value = 11

$ws.Range("A134").Value = "Mauricio"
$ws.Range("B134").Value = 332
$ws.Range("C134").Value = 0
$ws.Range("D134").Value = 322
$ws.Range("E134").Value = 0
$ws.Range("F134").Value = 0
$ws.Range("G134").Value = 0
$ws.Range("H134").Value = 10

$ws.Range("A135").Value = "Isla de Man"
$ws.Range("B135").Value = 330
$ws.Range("C135").Value = 0
$ws.Range("D135").Value = 271
$ws.Range("E135").Value = 36
$ws.Range("F135").Value = 21
$ws.Range("G135").Value = 0
$ws.Range("H135").Value = 23

$ws.Range("A136").Value = "Benin"
$ws.Range("B136").Value = 327
$ws.Range("C136").Value = 8
$ws.Range("D136").Value = 76
$ws.Range("E136").Value = 249
$ws.Range("F136").Value = 0
$ws.Range("G136").Value = 0
$ws.Range("H136").Value = 2

$ws.Range("A137").Value = "Montenegro"
$ws.Range("B137").Value = 324
$ws.Range("C137").Value = 0
$ws.Range("D137").Value = 298
$ws.Range("E137").Value = 17
$ws.Range("F137").Value = 2
$ws.Range("G137").Value = 0
$ws.Range("H137").Value = 9

# Rows 142-145: "Nepal" jumps ahead of Liberia/Haiti/Santo Tome y Principe,
# pushing each of those down one rank.
$ws.Range("A142").Value = "Nepal"
$ws.Range("B142").Value = 217
$ws.Range("C142").Value = 83
$ws.Range("D142").Value = 33
$ws.Range("E142").Value = 184
$ws.Range("F142").Value = 0
$ws.Range("G142").Value = 0
$ws.Range("H142").Value = 0

$ws.Range("A143").Value = "Liberia"
$ws.Range("B143").Value = 211
$ws.Range("C143").Value = 0
$ws.Range("D143").Value = 85
$ws.Range("E143").Value = 106
$ws.Range("F143").Value = 0
$ws.Range("G143").Value = 0
$ws.Range("H143").Value = 20

$ws.Range("A144").Value = "Haiti"
$ws.Range("B144").Value = 209
$ws.Range("C144").Value = 27
$ws.Range("D144").Value = 17
$ws.Range("E144").Value = 176
$ws.Range("F144").Value = 0
$ws.Range("G144").Value = 1
$ws.Range("H144").Value = 16

$ws.Range("A145").Value = "Santo Tome y Principe"
$ws.Range("B145").Value = 208
$ws.Range("C145").Value = 0
$ws.Range("D145").Value = 4
$ws.Range("E145").Value = 199
$ws.Range("F145").Value = 0
$ws.Range("G145").Value = 0
$ws.Range("H145").Value = 5
